$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 - VOLTAREN 75MG/3ML 3 AMP. (item 19)
$ws.Range("H25").Value = "5:2"

# P25 cell is formatted as a number (0.00) but must keep storing text,
# so force the format to text, write the value, then restore the format.
$fmt = $ws.Range("P25").NumberFormat
$ws.Range("P25").NumberFormat = "@"
$ws.Range("P25").Value = "50.4900"
$ws.Range("P25").NumberFormat = $fmt

$ws.Range("Q25").Value = "0:3"

# Row 31 - سرنجات 3 سم (item 25)
$fmt = $ws.Range("P31").NumberFormat
$ws.Range("P31").NumberFormat = "@"
$ws.Range("P31").Value = "15.6000"
$ws.Range("P31").NumberFormat = $fmt

$ws.Range("Q31").Value = "7:1"

# Row 34 - totals (recomputed sum after the price updates above)
$ws.Range("P34").Value = 1107.77

# Row 35 - footer timestamp, re-uploaded a few minutes later
$ws.Range("A35").Value = "Saturday, 7 June, 2025 7:31 PM"
